# Add a new stadium entry (VfR Krefeld) and correct the name_short of the
# existing "S.C. Viktoria 09 e.V." entry, per commit "add match protocol download".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stadiums")

# --- Correct name_short for the existing row 38 ("S.C. Viktoria 09 e.V.") ---
# It previously (incorrectly) held "Krefeld"; rename to "Vikt. Krefeld".
$ws.Range("C38").Value = "Vikt. Krefeld"

# --- Append a new row (39) for "VfR Krefeld 1920 e.V." ---
# Copy the formatting (styles, wrap text, date format, etc.) from the prior
# last row down into the new row before filling in its values.
$ws.Range("A38:J38").Copy()
$ws.Range("A39:J39").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "VfR Krefeld 1920 e.V."
$ws.Range("C39").Value = "VfR Krefeld"
$ws.Range("D39").Value = "https://goo.gl/maps/nFYDPTdWDCT2"
$ws.Range("E39").Value = "NULL"
$ws.Range("F39").Value = "NULL"
$ws.Range("G39").Value = "NULL"
$ws.Range("H39").Value = 1
$ws.Range("I39").Formula = "=NOW()"
$ws.Range("J39").Formula = "=NOW()"

# Match the author's final selection/scroll state.
$null = $ws.Range("E39:J39").Select()
